# Add two new rows documenting the Landlord possession claims dataset
# (MoJ Mortgage and landlord possession statistics, 2024 and 2021 editions)
# to the dataset guide sheet, consistent with the existing rows' layout:
# A=y, B=Variable, C=Table, D=Format, E=Years, F=Accessed, G=Link

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$accessedDate = Get-Date -Year 2025 -Month 9 -Day 2

# Row 27: 2024 possession statistics
$ws.Range("A27").Value = "y"
$ws.Range("B27").Value = "Landlord possession claims"
$ws.Range("C27").Value = "possession_statistics_2024"
$ws.Range("D27").Value = "MAP_CSV"
$ws.Range("E27").Value = 2024
$ws.Range("F27").Value = $accessedDate
$ws.Range("G27").Value = "https://assets.publishing.service.gov.uk/media/66b360d9a3c2a28abb50de35/Mortgage_and_landlord_statistical_data.zip"

# Row 28: 2021 possession statistics
$ws.Range("A28").Value = "y"
$ws.Range("B28").Value = "Landlord possession claims"
$ws.Range("C28").Value = "possession_statistics_2021"
$ws.Range("D28").Value = "MAP_CSV"
$ws.Range("E28").Value = 2021
$ws.Range("F28").Value = $accessedDate
$ws.Range("G28").Value = "https://assets.publishing.service.gov.uk/media/620397bce90e077f71cd545b/Mortgage_and_landlord_statistical_data.zip"

# Match existing formatting: column E left aligned (style index 2 used throughout),
# column F formatted as a date (style index 1 used throughout)
$ws.Range("E27:E28").Style = $ws.Range("E26").Style
$ws.Range("F27:F28").NumberFormat = $ws.Range("F26").NumberFormat

$ws.Range("A1:G28").Select()
$ws.Range("E29").Select()
